# Append three identical data rows (rows 2-4) below the existing header
# row on Sheet1, then leave the selection positioned the way the author
# left it (the whole of row 4 selected, active cell A4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title = "title"
$type  = "audio"
$owner = "مشارى"
$occasion = "حرب اكتوبر"
$path  = "1.mp3"

for ($r = 2; $r -le 4; $r++) {
    $ws.Range("A$r").Value = $title
    $ws.Range("B$r").Value = $type
    $ws.Range("C$r").Value = $owner
    $ws.Range("D$r").Value = $occasion
    $ws.Range("E$r").Value = $path
}

$ws.Range("A4:XFD4").Select()
